# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the regenerated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new value for column F, applied identically to both sheets.
$updates = @{
    3  = 3202
    5  = 130
    7  = 1688
    8  = 1633
    9  = 468
    10 = 370
    12 = 29
    17 = 232
    21 = 55
    24 = 214
    25 = 106
    26 = 33
    28 = 25
    29 = 266
    30 = 2178
    31 = 9
    34 = 330
    35 = 569
    40 = 518
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
